# Applies targeted cell updates (currentAveragePrice / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns, H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# as produced by the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 305
$ws.Range("I2").Value = 290
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 290
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -177
$ws.Range("N2").Value = -576
$ws.Range("H33").Value = 147.66667
$ws.Range("I33").Value = 122.25
$ws.Range("J33").Value = 198.5
$ws.Range("K33").Value = 122.25
$ws.Range("L33").Value = 198.5
$ws.Range("M33").Value = 106.75
$ws.Range("N33").Value = -656.5
$ws.Range("H55").Value = 420.13333
$ws.Range("I55").Value = 391.2
$ws.Range("K55").Value = 391.2
$ws.Range("M55").Value = -177.2
$ws.Range("H61").Value = 14689.6
$ws.Range("I61").Value = 69420
$ws.Range("J61").Value = 1007
$ws.Range("K61").Value = 208260
$ws.Range("L61").Value = 3021
$ws.Range("M61").Value = -208088
$ws.Range("N61").Value = -3365
$ws.Range("H64").Value = 10000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 10000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H82").Value = 665.5
$ws.Range("I82").Value = 665.5
$ws.Range("K82").Value = 1996.5
$ws.Range("M82").Value = -1590.5
$ws.Range("H85").Value = 665.5
$ws.Range("I85").Value = 665.5
$ws.Range("K85").Value = 1996.5
$ws.Range("M85").Value = -592.5
$ws.Range("H86").Value = 2667.6667
$ws.Range("I86").Value = 2003
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2003
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -880
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 2667.6667
$ws.Range("I89").Value = 2003
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 10015
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -4399
$ws.Range("N89").Value = -26232
$ws.Range("H100").Value = 1511.8334
$ws.Range("I100").Value = 1512.909
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1512.909
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -971.9090000000001
$ws.Range("N100").Value = -2582
$ws.Range("H132").Value = 12899.392
$ws.Range("I132").Value = 11409.789
$ws.Range("J132").Value = 19975
$ws.Range("K132").Value = 34229.367
$ws.Range("L132").Value = 59925
$ws.Range("M132").Value = -31699.367
$ws.Range("N132").Value = -64985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 5821
$ws.Range("I39").Value = 4985.2
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 4985.2
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -4465.2
$ws.Range("N39").Value = -11040
$ws.Range("H61").Value = 7499.25
$ws.Range("I61").Value = 7499.25
$ws.Range("K61").Value = 7499.25
$ws.Range("M61").Value = -7287.25
$ws.Range("H63").Value = 11831.75
$ws.Range("J63").Value = 17589.8
$ws.Range("L63").Value = 17589.8
$ws.Range("N63").Value = -18961.8
$ws.Range("H66").Value = 11831.75
$ws.Range("J66").Value = 17589.8
$ws.Range("L66").Value = 87949
$ws.Range("N66").Value = -94813
$ws.Range("H74").Value = 3062
$ws.Range("I74").Value = 3062
$ws.Range("K74").Value = 3062
$ws.Range("M74").Value = -2188
$ws.Range("H77").Value = 3062
$ws.Range("I77").Value = 3062
$ws.Range("K77").Value = 15310
$ws.Range("M77").Value = -10942
$ws.Range("H132").Value = 3681
$ws.Range("I132").Value = 3655.1936
$ws.Range("K132").Value = 10965.5808
$ws.Range("M132").Value = -8435.5808
$ws.Range("H136").Value = 7499.25
$ws.Range("I136").Value = 7499.25
$ws.Range("K136").Value = 22497.75
$ws.Range("M136").Value = -19947.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 9687.799999999999
$ws.Range("I97").Value = 9687.799999999999
$ws.Range("K97").Value = 9687.799999999999
$ws.Range("M97").Value = -8696.799999999999
$ws.Range("H134").Value = 2982.6667
$ws.Range("I134").Value = 2982.6667
$ws.Range("K134").Value = 8948.000100000001
$ws.Range("M134").Value = -6413.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3818.6743
$ws.Range("J31").Value = 7762.0713
$ws.Range("L31").Value = 7762.0713
$ws.Range("N31").Value = -8352.0713
$ws.Range("H34").Value = 3818.6743
$ws.Range("J34").Value = 7762.0713
$ws.Range("L34").Value = 7762.0713
$ws.Range("N34").Value = -8166.0713
$ws.Range("H58").Value = 3757.6155
$ws.Range("J58").Value = 4120
$ws.Range("L58").Value = 4120
$ws.Range("N58").Value = -4526
$ws.Range("H62").Value = 4677.5713
$ws.Range("J62").Value = 3435.75
$ws.Range("L62").Value = 3435.75
$ws.Range("N62").Value = -4683.75
$ws.Range("H65").Value = 4677.5713
$ws.Range("J65").Value = 3435.75
$ws.Range("L65").Value = 17178.75
$ws.Range("N65").Value = -23418.75
$ws.Range("H122").Value = 2424.75
$ws.Range("I122").Value = 2424.75
$ws.Range("K122").Value = 7274.25
$ws.Range("M122").Value = -4824.25
$ws.Range("H136").Value = 3757.6155
$ws.Range("J136").Value = 4120
$ws.Range("L136").Value = 12360
$ws.Range("N136").Value = -17460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1586.1818
$ws.Range("J5").Value = 1658
$ws.Range("L5").Value = 4974
$ws.Range("N5").Value = -5198
$ws.Range("H15").Value = 242.66667
$ws.Range("I15").Value = 310.66666
$ws.Range("K15").Value = 931.9999799999999
$ws.Range("M15").Value = -791.9999799999999
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H86").Value = 482.83334
$ws.Range("I86").Value = 482.83334
$ws.Range("K86").Value = 1448.50002
$ws.Range("M86").Value = -262.5000199999999
$ws.Range("H89").Value = 482.83334
$ws.Range("I89").Value = 482.83334
$ws.Range("K89").Value = 4345.50006
$ws.Range("M89").Value = 1582.49994
$ws.Range("H135").Value = 1586.1818
$ws.Range("J135").Value = 1658
$ws.Range("L135").Value = 14922
$ws.Range("N135").Value = -19992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3160.8333
$ws.Range("J80").Value = 2993.75
$ws.Range("L80").Value = 2993.75
$ws.Range("N80").Value = -4989.75
$ws.Range("H83").Value = 3160.8333
$ws.Range("J83").Value = 2993.75
$ws.Range("L83").Value = 14968.75
$ws.Range("N83").Value = -24952.75
$ws.Range("H102").Value = 1986.6786
$ws.Range("I102").Value = 1875.0741
$ws.Range("K102").Value = 1875.0741
$ws.Range("M102").Value = -253.0741
$ws.Range("H132").Value = 3238.6667
$ws.Range("I132").Value = 2687.5
$ws.Range("J132").Value = 5994.5
$ws.Range("K132").Value = 8062.5
$ws.Range("L132").Value = 17983.5
$ws.Range("M132").Value = -5532.5
$ws.Range("N132").Value = -23043.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 336.25
$ws.Range("I22").Value = 397.5
$ws.Range("J22").Value = 275
$ws.Range("K22").Value = 397.5
$ws.Range("L22").Value = 275
$ws.Range("M22").Value = -102.5
$ws.Range("N22").Value = -865
$ws.Range("H27").Value = 336.25
$ws.Range("I27").Value = 397.5
$ws.Range("J27").Value = 275
$ws.Range("K27").Value = 397.5
$ws.Range("L27").Value = 275
$ws.Range("M27").Value = -290.5
$ws.Range("N27").Value = -489
$ws.Range("H46").Value = 4570
$ws.Range("I46").Value = 2633.3333
$ws.Range("J46").Value = 6022.5
$ws.Range("K46").Value = 2633.3333
$ws.Range("L46").Value = 6022.5
$ws.Range("M46").Value = -2445.3333
$ws.Range("N46").Value = -6398.5
$ws.Range("H68").Value = 6566.6665
$ws.Range("I68").Value = 2900
$ws.Range("K68").Value = 2900
$ws.Range("M68").Value = -2151
$ws.Range("H71").Value = 6566.6665
$ws.Range("I71").Value = 2900
$ws.Range("K71").Value = 14500
$ws.Range("M71").Value = -10756

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7500
$ws.Range("J81").Value = 7500
$ws.Range("L81").Value = 15000
$ws.Range("N81").Value = -17122
$ws.Range("H84").Value = 7500
$ws.Range("J84").Value = 7500
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -85608
$ws.Range("H96").Value = 706.6667
$ws.Range("I96").Value = 449.2
$ws.Range("K96").Value = 449.2
$ws.Range("M96").Value = 923.8
$ws.Range("H100").Value = 725.8
$ws.Range("I100").Value = 473.1111
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 946.2222
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -405.2222
$ws.Range("N100").Value = -7082

